$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I and J), matching the formatting of
# the existing header cells in row 1 (bold font, thin border, centered /
# top-aligned), mirroring H1's formatting.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

foreach ($addr in @("I1", "J1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data values for columns I and J, rows 2-69.
$data = @(
    @{Row=2; I=8; J=8},
    @{Row=3; I=8; J=8},
    @{Row=4; I=9; J=9},
    @{Row=5; I=8; J=8},
    @{Row=6; I=9; J=9},
    @{Row=7; I=6; J=7},
    @{Row=8; I=7; J=7},
    @{Row=9; I=8; J=8},
    @{Row=10; I=8; J=8},
    @{Row=11; I=6; J=7},
    @{Row=12; I=8; J=8},
    @{Row=13; I=8; J=8},
    @{Row=14; I=9; J=9},
    @{Row=15; I=9; J=9},
    @{Row=16; I=8; J=8},
    @{Row=17; I=7; J=7},
    @{Row=18; I=9; J=9},
    @{Row=19; I=8; J=8},
    @{Row=20; I=7; J=8},
    @{Row=21; I=10; J=10},
    @{Row=22; I=9; J=9},
    @{Row=23; I=9; J=9},
    @{Row=24; I=9; J=9},
    @{Row=25; I=7; J=7},
    @{Row=26; I=8; J=8},
    @{Row=27; I=7; J=7},
    @{Row=28; I=8; J=8},
    @{Row=29; I=7; J=7},
    @{Row=30; I=7; J=7},
    @{Row=31; I=7; J=7},
    @{Row=32; I=7; J=8},
    @{Row=33; I=7; J=7},
    @{Row=34; I=7; J=8},
    @{Row=35; I=8; J=8},
    @{Row=36; I=7; J=7},
    @{Row=37; I=7; J=7},
    @{Row=38; I=7; J=7},
    @{Row=39; I=8; J=9},
    @{Row=40; I=7; J=7},
    @{Row=41; I=8; J=8},
    @{Row=42; I=6; J=7},
    @{Row=43; I=7; J=8},
    @{Row=44; I=6; J=6},
    @{Row=45; I=9; J=9},
    @{Row=46; I=7; J=8},
    @{Row=47; I=8; J=8},
    @{Row=48; I=9; J=9},
    @{Row=49; I=7; J=7},
    @{Row=50; I=7; J=7},
    @{Row=51; I=7; J=8},
    @{Row=52; I=6; J=6},
    @{Row=53; I=8; J=8},
    @{Row=54; I=7; J=7},
    @{Row=55; I=7; J=7},
    @{Row=56; I=8; J=8},
    @{Row=57; I=6; J=6},
    @{Row=58; I=7; J=8},
    @{Row=59; I=7; J=7},
    @{Row=60; I=6; J=6},
    @{Row=61; I=8; J=8},
    @{Row=62; I=7; J=7},
    @{Row=63; I=8; J=8},
    @{Row=64; I=7; J=7},
    @{Row=65; I=7; J=7},
    @{Row=66; I=8; J=8},
    @{Row=67; I=4; J=4},
    @{Row=68; I=3; J=3},
    @{Row=69; I=3; J=3}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
